# The deck ships two embedded theme parts:
#   ppt/theme/theme2.xml  -> bound to the slide master (this is the theme
#                             that is actually "applied" to the deck) and
#                             currently holds the "Integral" colour scheme.
#   ppt/theme/theme1.xml  -> bound to the notes master, holding the
#                             "Office Theme" colour scheme.
#
# The author re-themed the presentation (Design tab -> "Office Theme"),
# which swapped the two colour palettes between the two theme parts.
# Re-create that by pushing the standard Office theme colours onto the
# presentation's live theme colour scheme -- every accessor (slide
# master / notes master / handout master / individual slides) shares
# this single colour scheme object, so setting it once re-themes the
# whole deck.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office Theme palette (RGB hex, written here as BGR COM colour literals):
#   dk1      000000
#   lt1      FFFFFF
#   dk2      44546A
#   lt2      E7E6E6
#   accent1  5B9BD5
#   accent2  ED7D31
#   accent3  A5A5A5
#   accent4  FFC000
#   accent5  4472C4
#   accent6  70AD47
#   hlink    0563C1
#   folHlink 954F72
$colors.Item(1).RGB = 0x000000
$colors.Item(2).RGB = 0xFFFFFF
$colors.Item(3).RGB = 0x6A5444
$colors.Item(4).RGB = 0xE6E6E7
$colors.Item(5).RGB = 0xD59B5B
$colors.Item(6).RGB = 0x317DED
$colors.Item(7).RGB = 0xA5A5A5
$colors.Item(8).RGB = 0x00C0FF
$colors.Item(9).RGB = 0xC47244
$colors.Item(10).RGB = 0x47AD70
$colors.Item(11).RGB = 0xC16305
$colors.Item(12).RGB = 0x724F95

# Best-effort: update the in-memory theme/design name too (some hosts
# only persist the colour values, but this keeps object-model reads
# consistent within the session).
$theme.Name = "Office Theme"
